# The deck ships with two embedded themes: the applied design theme
# ("Integral", theme1.xml, used by the Slide Master) and the default
# "Office Theme" palette (theme2.xml, used only by the Notes Master).
#
# This edit switches the presentation's design theme over to that
# "Office Theme" colour palette - i.e. every one of the 12 theme colour
# slots (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink) is repointed at the
# Office Theme's RGB values, the same values already sitting in the
# Notes Master's theme.

$p  = $ppt.ActivePresentation
$d  = $p.Designs.Item(1)
$cs = $d.SlideMaster.Theme.ThemeColorScheme

# Office Theme palette, in ThemeColorScheme.Colors(1..12) slot order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
$officeTheme = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

for ($i = 0; $i -lt $officeTheme.Count; $i++) {
    $hex = $officeTheme[$i]
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)

    # OLE RGB() packing: 0x00BBGGRR
    $cs.Colors($i + 1).RGB = $r + ($g * 256) + ($b * 65536)
}
